# Update Montecreto sheet with new daily rows (data through 2022-01-05),
# matching commit "aggiornamento fino a 6 gennaio 2022".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows: serial date (col A), nuovi pos. (col B),
# somma mobile 7gg. (col C), somma mobile 7gg. per 100mila abitanti (col D)
$newRows = @(
    @(44539, 0, 0, 0),
    @(44540, 0, 0, 0),
    @(44541, 0, 0, 0),
    @(44542, 0, 0, 0),
    @(44543, 0, 0, 0),
    @(44544, 1, 1, 109.1703056768559),
    @(44545, 0, 1, 109.1703056768559),
    @(44546, 2, 3, 327.5109170305677),
    @(44547, 0, 3, 327.5109170305677),
    @(44548, 1, 4, 436.6812227074236),
    @(44550, 1, 5, 545.8515283842795),
    @(44551, 1, 6, 655.0218340611353),
    @(44552, 0, 5, 545.8515283842795),
    @(44553, 0, 5, 545.8515283842795),
    @(44554, 1, 4, 436.6812227074236),
    @(44555, 1, 5, 545.8515283842795),
    @(44556, 2, 6, 655.0218340611353),
    @(44557, 0, 5, 545.8515283842795),
    @(44558, 0, 4, 436.6812227074236),
    @(44559, 2, 6, 655.0218340611353),
    @(44560, 0, 6, 655.0218340611353),
    @(44561, 0, 5, 545.8515283842795),
    @(44562, 1, 5, 545.8515283842795),
    @(44563, 0, 3, 327.5109170305677),
    @(44564, 0, 3, 327.5109170305677),
    @(44565, 0, 3, 327.5109170305677),
    @(44566, 0, 1, 109.1703056768559)
)

$startRow = 465
$endRow = $startRow + $newRows.Count - 1

# Extend column A formatting (date style used throughout the column) down
# over the new rows before writing values, so the new cells inherit the
# same style as the rest of the sheet rather than getting a fresh one.
$ws.Range("A464").Copy($ws.Range("A" + $startRow + ":A" + $endRow))

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $vals = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
}
